# Lessee Information.xlsx edit script
# Commit message: "Added Units, Fixed UI notifications not displaying properly"
#
# Summary of data changes (derived from the OOXML diff):
#   - A1 (unit "1A") was a completely blank row; B/C/D now read "Null "
#   - Unit 2A (was "Null "/"Null "/"Null ") now leased to Cris Dione Sigua / Residential / Null
#   - Unit 3E (was "Null "/"Null "/"Null ") now leased to Carissa Tapang / Residential / Null
#   - Unit 3G (was "Null "/"Null "/"Null ") now leased to Misie Quimba / Residential / Null
#   - Unit 4C (was "Vacant") now leased to Marvin Inocencio
#   - Unit 4H (was "Null "/"Null "/"Null ") now leased to Jonald Cajilig / Residential / Null
#   - 12 new units added: 2I,2J,2K,2L (after 2H); 3I,3J,3K,3L (after 3H); 4I,4J,4K,4L (appended at end)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in unit 1A, which previously had completely empty cells ---
$ws.Range("B2").Value = "Null "
$ws.Range("C2").Value = "Null "
$ws.Range("D2").Value = "Null "

# --- Previously-vacant units now have lessees ---
$ws.Range("B6").Value = "Cris Dione Sigua"
$ws.Range("C6").Value = "Residential"

$ws.Range("B18").Value = "Carissa Tapang"
$ws.Range("C18").Value = "Residential"

$ws.Range("B20").Value = "Misie Quimba"
$ws.Range("C20").Value = "Residential"

$ws.Range("B24").Value = "Marvin Inocencio"

$ws.Range("B29").Value = "Jonald Cajilig"
$ws.Range("C29").Value = "Residential"

# --- Insert 4 new units (2I, 2J, 2K, 2L) right after 2H (row 13) ---
$ws.Rows.Item(14).Resize(4).Insert()
$ws.Range("A14").Value = "2I"
$ws.Range("B14").Value = "Angie Villarico"
$ws.Range("C14").Value = "Residential"
$ws.Range("D14").Value = "Null "

$ws.Range("A15").Value = "2J"
$ws.Range("B15").Value = "Ronwaldo Bariuan"
$ws.Range("C15").Value = "Residential"
$ws.Range("D15").Value = "Null "

$ws.Range("A16").Value = "2K"
$ws.Range("B16").Value = "Cecille Espiritu"
$ws.Range("C16").Value = "Residential"
$ws.Range("D16").Value = "Null "

$ws.Range("A17").Value = "2L"
$ws.Range("B17").Value = "Arvin de Guzman"
$ws.Range("C17").Value = "Residential"
$ws.Range("D17").Value = "Null "

# --- Insert 4 new units (3I, 3J, 3K, 3L) right after 3H (now row 25) ---
$ws.Rows.Item(26).Resize(4).Insert()
$ws.Range("A26").Value = "3I"
$ws.Range("B26").Value = "Jennifer Valenzuela"
$ws.Range("C26").Value = "Residential"
$ws.Range("D26").Value = "Null "

$ws.Range("A27").Value = "3J"
$ws.Range("B27").Value = "Clarissa Gallardo"
$ws.Range("C27").Value = "Residential"
$ws.Range("D27").Value = "Null "

$ws.Range("A28").Value = "3K"
$ws.Range("B28").Value = "Aaron Alfonso"
$ws.Range("C28").Value = "Residential"
$ws.Range("D28").Value = "Null "

$ws.Range("A29").Value = "3L"
$ws.Range("B29").Value = "Rowena Barcelona"
$ws.Range("C29").Value = "Residential"
$ws.Range("D29").Value = "Null "

# --- Append 4 new units (4I, 4J, 4K, 4L) at the very end (now row 37 is last, 4H) ---
$ws.Range("A38").Value = "4I"
$ws.Range("B38").Value = "Rhodora Capulong"
$ws.Range("C38").Value = "Residential"
$ws.Range("D38").Value = "Null "

$ws.Range("A39").Value = "4J"
$ws.Range("B39").Value = "Philip Sevilla"
$ws.Range("C39").Value = "Residential"
$ws.Range("D39").Value = "Null "

$ws.Range("A40").Value = "4K"
$ws.Range("B40").Value = "Karren Ralutin"
$ws.Range("C40").Value = "Residential"
$ws.Range("D40").Value = "Null "

$ws.Range("A41").Value = "4L"
$ws.Range("B41").Value = "Zaida Rosales"
$ws.Range("C41").Value = "Residential"
$ws.Range("D41").Value = "Null "

# --- Update the current selection to match the authored file (cosmetic) ---
$ws.Range("E20").Select()
